$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect new date
$ws.Name = "Through 2021-11-25"

# Update the "November (through 11-24)" label to 11-25
$ws.Range("A12").Value = "November (through 11-25)"

# Update November row (row 12) values for columns C:H (2016-2021)
$ws.Range("C12").Value = 62
$ws.Range("D12").Value = 94
$ws.Range("E12").Value = 50
$ws.Range("F12").Value = 45
$ws.Range("G12").Value = 180
$ws.Range("H12").Value = 171

# Update Total row (row 13) values for columns C:H (2016-2021)
$ws.Range("C13").Value = 548
$ws.Range("D13").Value = 804
$ws.Range("E13").Value = 665
$ws.Range("F13").Value = 527
$ws.Range("G13").Value = 1237
$ws.Range("H13").Value = 1614
